# Generate Report for Handback
# Updates the timestamp cells on the Overview / zh-cn / de-de sheets to
# reflect the new handback/handoff generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$wsOverview.Range("G2").Value = "2017-02-21 10:37:05"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn.Range("H2").Value = "2017-02-21 10:36:47"
$wsZhCn.Range("L2").Value = "2017-02-21 10:37:46"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsDeDe.Range("H2").Value = "2017-02-21 10:37:05"
$wsDeDe.Range("L2").Value = "2017-02-21 10:38:09"
